$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.359.73'
$ws.Range("E2").Value = '  -1.25%  '

$ws.Range("D3").Value = '1.798.80'
$ws.Range("E3").Value = '  -1.45%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '227.02'
$ws.Range("E5").Value = '  -0.88%  '

$ws.Range("E6").Value = '  +3.42%  '

$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '36.23'
$ws.Range("E8").Value = '  +3.76%  '

$ws.Range("E9").Value = '  -2.19%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.0688'
$ws.Range("E10").Value = '  -1.72%  '

$ws.Range("E11").Value = '  +1.06%  '

$ws.Range("D12").Value = '2.057.82'
$ws.Range("E12").Value = '  -1.50%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '11.28'
$ws.Range("E13").Value = '  -1.67%  '

$ws.Range("D14").Value = '1.826.18'
$ws.Range("E14").Value = '  +0.43%  '

$ws.Range("E15").Value = '  -1.20%  '

$ws.Range("D16").Value = '34.347.24'
$ws.Range("E16").Value = '  -1.19%  '

$ws.Range("E17").Value = '  +1.67%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '69.69'
$ws.Range("E18").Value = '  +0.80%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '243.81'
$ws.Range("E19").Value = '  -1.65%  '

$ws.Range("D20").Value = '0.0₃0783'
$ws.Range("E20").Value = '  -2.44%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '11.41'
$ws.Range("E21").Value = '  -2.06%  '

$ws.Range("E22").Value = '  +0.33%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '4.15'
$ws.Range("E23").Value = '  -1.27%  '

$ws.Range("E24").Value = '  +6.26%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '170.22'
$ws.Range("E25").Value = '  -2.23%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '8.08'
$ws.Range("E26").Value = '  +7.74%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '17.25'
$ws.Range("E27").Value = '  +2.15%  '

$ws.Range("E28").Value = '  +1.59%  '

$ws.Range("E29").Value = '  +0.26%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '3.94'
$ws.Range("E30").Value = '  -1.93%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '1.24'
$ws.Range("E31").Value = '  -1.24%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '3.81'
$ws.Range("E32").Value = '  -1.49%  '

$ws.Range("E33").Value = '  -2.92%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '1.79'
$ws.Range("E34").Value = '  -4.32%  '

$ws.Range("D35").Value = '1.373.53'
$ws.Range("E35").Value = '  -3.09%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.655'
$ws.Range("E36").Value = '  -4.68%  '

$ws.Range("E37").Value = '  -2.41%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '2.36'
$ws.Range("E38").Value = '  -10.78%  '

$ws.Range("E39").Value = '  -3.33%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '81.76'
$ws.Range("E40").Value = '  -4.49%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '2.41'
$ws.Range("E41").Value = '  +0.52%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.945'
$ws.Range("E42").Value = '  -1.77%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '2.79'
$ws.Range("E43").Value = '  -2.79%  '

$ws.Range("E44").Value = '  +6.04%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '13.30'
$ws.Range("E45").Value = '  -4.07%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.0496'
$ws.Range("E46").Value = '  -4.26%  '

$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.959.18'
$ws.Range("E47").Value = '  -1.48%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '5.84'
$ws.Range("E48").Value = '  -4.44%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '102.66'
$ws.Range("E50").Value = '  -3.01%  '

$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '48.79'
$ws.Range("E51").Value = '  -3.03%  '
